# Auto-generated edit script: updates cryptocurrency Price (D) and
# Volume(1h) (E) columns for rows 2-51 per the scraped GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "60.151.71"
$ws.Range("E2").Value = "  -5.59%  "
$ws.Range("D3").Value = "3.341.96"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "564.97"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").Value = "130.57"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.340.72"
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").Value = "3.914.20"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "3.344.24"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").Value = "0.0000169"
$ws.Range("E16").Value = "  -3.62%  "
$ws.Range("D17").Value = "24.57"
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "60.240.22"
$ws.Range("E18").Value = "  -5.42%  "
$ws.Range("D19").Value = "5.69"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "13.45"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  -7.28%  "
$ws.Range("D22").Value = "354.67"
$ws.Range("E22").Value = "  -7.76%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "3.476.06"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "69.39"
$ws.Range("E26").Value = "  -6.27%  "
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").Value = "1.68"
$ws.Range("E28").Value = "  +19.07%  "
$ws.Range("D29").Value = "7.51"
$ws.Range("E29").Value = "  +7.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "7.95"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").Value = "0.154"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "3.373.75"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").Value = "22.98"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "5.42"
$ws.Range("E37").Value = "  +5.91%  "
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").Value = "159.15"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("D41").Value = "0.0769"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").Value = "  +9.17%  "
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").Value = "40.78"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "23.79"
$ws.Range("E47").Value = "  +2.30%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").Value = "22.41"
$ws.Range("E50").Value = "  +10.42%  "
$ws.Range("D51").Value = "0.893"
$ws.Range("E51").Value = "  +0.60%  "
